# Y4_B2526_Excuses.xlsx update
# - Replace Student ID values (column A) for existing rows 2-23 with new
#   values pulled in by the attendance app.
# - Append 4 new "general surgery" Excuse log rows (24-27) with the same
#   layout/style as existing rows, extending the sheet dimension to F27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Student ID values (column A) ---
$ws.Cells.Item(2,1).Value  = "190981"
$ws.Cells.Item(3,1).Value  = "191055"
$ws.Cells.Item(4,1).Value  = "211216"
$ws.Cells.Item(5,1).Value  = "201669"
$ws.Cells.Item(6,1).Value  = "190922"
$ws.Cells.Item(7,1).Value  = "211137"
$ws.Cells.Item(8,1).Value  = "200785"
$ws.Cells.Item(9,1).Value  = "200116"
$ws.Cells.Item(10,1).Value = "201632"
$ws.Cells.Item(11,1).Value = "201563"
$ws.Cells.Item(12,1).Value = "180804"
$ws.Cells.Item(13,1).Value = "190807"
$ws.Cells.Item(14,1).Value = "191109"
$ws.Cells.Item(15,1).Value = "210923"
$ws.Cells.Item(16,1).Value = "201026"
$ws.Cells.Item(17,1).Value = "181013"
$ws.Cells.Item(18,1).Value = "201157"
$ws.Cells.Item(19,1).Value = "211096"
$ws.Cells.Item(20,1).Value = "211147"
$ws.Cells.Item(21,1).Value = "211046"
$ws.Cells.Item(22,1).Value = "190803"
$ws.Cells.Item(23,1).Value = "201572"

# --- Append new rows 24-27, copying formats (fill/font/alignment) from the
#     last two existing data rows (22 = even style, 23 = odd style) so the
#     alternating row styling continues correctly ---
$ws.Range("A22:F22").Copy()
$ws.Range("A24:F24").PasteSpecial(-4122)

$ws.Range("A23:F23").Copy()
$ws.Range("A25:F25").PasteSpecial(-4122)

$ws.Range("A22:F22").Copy()
$ws.Range("A26:F26").PasteSpecial(-4122)

$ws.Range("A23:F23").Copy()
$ws.Range("A27:F27").PasteSpecial(-4122)

$newRows = @(
    @{Row=24; Id="211133"},
    @{Row=25; Id="200228"},
    @{Row=26; Id="191131"},
    @{Row=27; Id="200869"}
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r,1).Value = $item.Id
    $ws.Cells.Item($r,2).Value = "general surgery"
    $ws.Cells.Item($r,3).Value = "25/10/2025"
    $ws.Cells.Item($r,4).Value = "10:30:00"
    $ws.Cells.Item($r,5).Value = "Excuse"
    $ws.Cells.Item($r,6).Value = "System"
}
